$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new VIN entry as the next row (row 28)
$ws.Range("A28").Value = "4V4MC9DG8DN542142"
$ws.Range("B28").Value = 100116

# Match formatting of the preceding data row (B27) for the new cell
$ws.Range("B27").Copy()
$ws.Range("B28").PasteSpecial(-4122)

# Update the active selection to match the post-edit workbook state
$ws.Range("C21").Select()
